# Implementing Select Date feature
# Move the Dep* config values (row 2, cols E:I) down to row 4, freeing
# E2:I2 for the new "Select Date" columns, and update the active
# selection/view to point at the relocated range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing values that live in E2:I2 before clearing them.
$depFile = $ws.Range("E2").Value2
$depSheet = $ws.Range("F2").Value2
$depColCount = $ws.Range("G2").Value2
$depKeyCol = $ws.Range("H2").Value2
$depValCol = $ws.Range("I2").Value2

# Clear the old location.
$ws.Range("E2:I2").ClearContents()

# Re-home the values on row 4.
$ws.Range("E4").Value = $depFile
$ws.Range("F4").Value = $depSheet
$ws.Range("G4").Value = $depColCount
$ws.Range("H4").Value = $depKeyCol
$ws.Range("I4").Value = $depValCol

# Update the view: scroll so column C is the left-most visible column,
# and select the relocated range E4:I4 (active cell E4).
$ws.Range("E4:I4").Select()
$ws.Application.ActiveWindow.ScrollColumn = 3
